# Add Armor Potion script functionality:
# - Set the "scr" column (C) for the armor potion rows (6-9) to "scr_armor_potion"
#   (previously a placeholder numeric value of -1, meaning "no script")
# - Move the active cell selection from C13 to C14 (cosmetic selection state)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 6..9) {
    $ws.Cells.Item($row, 3).Value = "scr_armor_potion"
}

$ws.Range("C14").Select()
